# Auto-generated Excel COM-interop script to apply data updates
# to the Behemoth_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (18 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2441.6
$ws.Range("I70").Value = 1594.4
$ws.Range("J70").Value = 3288.8
$ws.Range("K70").Value = 4783.200000000001
$ws.Range("L70").Value = 9866.400000000001
$ws.Range("M70").Value = -4513.200000000001
$ws.Range("N70").Value = -10406.4
$ws.Range("H73").Value = 2441.6
$ws.Range("I73").Value = 1594.4
$ws.Range("J73").Value = 3288.8
$ws.Range("K73").Value = 4783.200000000001
$ws.Range("L73").Value = 9866.400000000001
$ws.Range("M73").Value = -3847.200000000001
$ws.Range("N73").Value = -11738.4
$ws.Range("H132").Value = 2306.85
$ws.Range("I132").Value = 1993.5714
$ws.Range("K132").Value = 5980.7142
$ws.Range("M132").Value = -3450.7142

# --- Sheet: ARM (49 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 30861.268
$ws.Range("J8").Value = 31260.572
$ws.Range("L8").Value = 31260.572
$ws.Range("N8").Value = -31548.572
$ws.Range("H32").Value = 15157961
$ws.Range("I32").Value = 15157961
$ws.Range("K32").Value = 15157961
$ws.Range("M32").Value = -15157674
$ws.Range("H61").Value = 20045734
$ws.Range("I61").Value = 33336880
$ws.Range("J61").Value = 109018.6
$ws.Range("K61").Value = 33336880
$ws.Range("L61").Value = 109018.6
$ws.Range("M61").Value = -33336668
$ws.Range("N61").Value = -109442.6
$ws.Range("H74").Value = 10007546
$ws.Range("I74").Value = 13889767
$ws.Range("K74").Value = 13889767
$ws.Range("M74").Value = -13888893
$ws.Range("H77").Value = 10007546
$ws.Range("I77").Value = 13889767
$ws.Range("K77").Value = 69448835
$ws.Range("M77").Value = -69444467
$ws.Range("H97").Value = 1443
$ws.Range("J97").Value = 1200
$ws.Range("L97").Value = 1200
$ws.Range("N97").Value = -2192
$ws.Range("H102").Value = 9748.666999999999
$ws.Range("I102").Value = 10798.5
$ws.Range("J102").Value = 4499.5
$ws.Range("K102").Value = 10798.5
$ws.Range("L102").Value = 4499.5
$ws.Range("M102").Value = -9176.5
$ws.Range("N102").Value = -7743.5
$ws.Range("H115").Value = 75051.5
$ws.Range("J115").Value = 75051.5
$ws.Range("L115").Value = 75051.5
$ws.Range("N115").Value = -78185.5
$ws.Range("H122").Value = 2224.125
$ws.Range("I122").Value = 1970.4286
$ws.Range("K122").Value = 5911.2858
$ws.Range("M122").Value = -3461.2858
$ws.Range("H136").Value = 20045734
$ws.Range("I136").Value = 33336880
$ws.Range("J136").Value = 109018.6
$ws.Range("K136").Value = 100010640
$ws.Range("L136").Value = 327055.8
$ws.Range("M136").Value = -100008090
$ws.Range("N136").Value = -332155.8

# --- Sheet: BSM (23 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 57494.5
$ws.Range("J27").Value = 57494.5
$ws.Range("L27").Value = 57494.5
$ws.Range("N27").Value = -57878.5
$ws.Range("H99").Value = 2262.6875
$ws.Range("I99").Value = 1807.6666
$ws.Range("K99").Value = 1807.6666
$ws.Range("M99").Value = -309.6666
$ws.Range("H105").Value = 1467.6666
$ws.Range("I105").Value = 1543.8572
$ws.Range("J105").Value = 1361
$ws.Range("K105").Value = 1543.8572
$ws.Range("L105").Value = 1361
$ws.Range("M105").Value = 203.1428000000001
$ws.Range("N105").Value = -4855
$ws.Range("H107").Value = 2219.0908
$ws.Range("I107").Value = 1925.25
$ws.Range("K107").Value = 1925.25
$ws.Range("M107").Value = -5.25
$ws.Range("H125").Value = 99323
$ws.Range("J125").Value = 99323
$ws.Range("L125").Value = 99323
$ws.Range("N125").Value = -109163

# --- Sheet: CRP (30 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1168546.8
$ws.Range("I31").Value = 1973
$ws.Range("K31").Value = 1973
$ws.Range("M31").Value = -1678
$ws.Range("H34").Value = 1168546.8
$ws.Range("I34").Value = 1973
$ws.Range("K34").Value = 1973
$ws.Range("M34").Value = -1771
$ws.Range("H107").Value = 1673.2307
$ws.Range("I107").Value = 1114.2858
$ws.Range("J107").Value = 2325.3333
$ws.Range("K107").Value = 1114.2858
$ws.Range("L107").Value = 2325.3333
$ws.Range("M107").Value = 805.7141999999999
$ws.Range("N107").Value = -6165.3333
$ws.Range("H108").Value = 75477.664
$ws.Range("J108").Value = 75477.664
$ws.Range("L108").Value = 75477.664
$ws.Range("N108").Value = -83157.664
$ws.Range("H117").Value = 49974
$ws.Range("J117").Value = 49974
$ws.Range("L117").Value = 49974
$ws.Range("N117").Value = -59152
$ws.Range("H132").Value = 1670.8
$ws.Range("I132").Value = 1670.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5012.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2482.4
$ws.Range("N132").ClearContents()

# --- Sheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1113046.2
$ws.Range("J92").Value = 2123.25
$ws.Range("L92").Value = 6369.75
$ws.Range("N92").Value = -8865.75
$ws.Range("H107").Value = 665.3333
$ws.Range("J107").Value = 786.2857
$ws.Range("L107").Value = 2358.8571
$ws.Range("N107").Value = -6198.8571
$ws.Range("H132").Value = 2379.2
$ws.Range("I132").Value = 2439.2666
$ws.Range("J132").Value = 2199
$ws.Range("K132").Value = 21953.3994
$ws.Range("L132").Value = 19791
$ws.Range("M132").Value = -19423.3994
$ws.Range("N132").Value = -24851

# --- Sheet: GSM (48 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 28027
$ws.Range("I80").Value = 24199.572
$ws.Range("J80").Value = 34725
$ws.Range("K80").Value = 24199.572
$ws.Range("L80").Value = 34725
$ws.Range("M80").Value = -23201.572
$ws.Range("N80").Value = -36721
$ws.Range("H83").Value = 28027
$ws.Range("I83").Value = 24199.572
$ws.Range("J83").Value = 34725
$ws.Range("K83").Value = 120997.86
$ws.Range("L83").Value = 173625
$ws.Range("M83").Value = -116005.86
$ws.Range("N83").Value = -183609
$ws.Range("H113").Value = 4022.75
$ws.Range("I113").Value = 3801.2307
$ws.Range("K113").Value = 3801.2307
$ws.Range("M113").Value = -1631.2307
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 49677.4
$ws.Range("J118").Value = 49677.4
$ws.Range("L118").Value = 49677.4
$ws.Range("N118").Value = -52991.4
$ws.Range("H122").Value = 2474.5
$ws.Range("I122").Value = 2474.5
$ws.Range("K122").Value = 7423.5
$ws.Range("M122").Value = -4973.5
$ws.Range("I126").Value = 2999
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8997
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6527
$ws.Range("N126").ClearContents()
$ws.Range("H128").Value = 120995
$ws.Range("J128").Value = 120995
$ws.Range("L128").Value = 120995
$ws.Range("N128").Value = -130955
$ws.Range("H129").Value = 69975
$ws.Range("J129").Value = 69975
$ws.Range("L129").Value = 69975
$ws.Range("N129").Value = -79975
$ws.Range("H132").Value = 55557988
$ws.Range("I132").Value = 58825956
$ws.Range("K132").Value = 176477868
$ws.Range("M132").Value = -176475338

# --- Sheet: LTW (23 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23823.53
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 25125
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 25125
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -25397
$ws.Range("H116").Value = 172518
$ws.Range("J116").Value = 172518
$ws.Range("L116").Value = 172518
$ws.Range("N116").Value = -181696
$ws.Range("H122").Value = 6844.2666
$ws.Range("I122").Value = 5962.6665
$ws.Range("K122").Value = 17887.9995
$ws.Range("M122").Value = -15437.9995
$ws.Range("H125").Value = 99896
$ws.Range("J125").Value = 99896
$ws.Range("L125").Value = 99896
$ws.Range("N125").Value = -109736
$ws.Range("H132").Value = 32744.777
$ws.Range("I132").Value = 5108.7715
$ws.Range("K132").Value = 15326.3145
$ws.Range("M132").Value = -12796.3145

# --- Sheet: WVR (23 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1549.25
$ws.Range("I113").Value = 1549.25
$ws.Range("K113").Value = 4647.75
$ws.Range("M113").Value = -2477.75
$ws.Range("H122").Value = 4293.8965
$ws.Range("I122").Value = 2964.6875
$ws.Range("K122").Value = 8894.0625
$ws.Range("M122").Value = -6444.0625
$ws.Range("H123").Value = 60214.5
$ws.Range("J123").Value = 60214.5
$ws.Range("L123").Value = 60214.5
$ws.Range("N123").Value = -70014.5
$ws.Range("H130").Value = 67122
$ws.Range("J130").Value = 69496
$ws.Range("L130").Value = 69496
$ws.Range("N130").Value = -79536
$ws.Range("H132").Value = 3850
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 3600
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -1070
$ws.Range("N132").Value = -24560

Write-Host "Applied updates to $($wb.Worksheets.Count) worksheets."